$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '62.262.95'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +1.46%  '

# Row 3
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.424.96'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +1.99%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '563.48'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +2.09%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '144.48'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +3.20%  '

# Row 7
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.533'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +1.76%  '

# Row 9
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '2.423.94'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +1.87%  '

# Row 10
$ws.Range('E10').Value = '  +2.01%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.154'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -2.16%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '5.38'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +0.36%  '

# Row 13
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.353'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +0.48%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '26.01'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +1.81%  '

# Row 15
$ws.Range('E15').Value = '  +5.69%  '

# Row 16
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '2.863.21'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +1.99%  '

# Row 17
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '62.189.01'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +1.34%  '

# Row 18
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '2.425.56'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +2.01%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '11.36'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +3.38%  '

# Row 20
$ws.Range('E20').Value = '  +1.09%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '324.68'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +1.10%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.77'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +0.92%  '

# Row 23
$ws.Range('E23').Value = '  -0.03%  '

# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '65.60'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +1.88%  '

# Row 25
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '1.72'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -2.43%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '8.92'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +0.66%  '

# Row 27
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '588.87'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +13.44%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0948'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +5.33%  '

# Row 29
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +0.07%  '

# Row 30
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '2.528.53'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +1.36%  '

# Row 31
$ws.Range('E31').Value = '  +5.71%  '

# Row 32
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '8.28'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +0.93%  '

# Row 33
$ws.Range('E33').Value = '  +0.37%  '

# Row 34
$ws.Range('E34').Value = '  +1.99%  '

# Row 35
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.57'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +1.84%  '

# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '5.76'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +4.44%  '

# Row 37
$ws.Range('E37').Value = '  +0.02%  '

# Row 38
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '4.82'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +2.56%  '

# Row 39
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.384'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +1.57%  '

# Row 40
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '153.78'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +4.66%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '18.70'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +0.91%  '

# Row 42
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.84'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -2.46%  '

# Row 43
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'

# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '2.35'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +8.40%  '

# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '150.43'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +1.47%  '

# Row 46
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '3.67'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +1.56%  '

# Row 47
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.0541'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +2.76%  '

# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '20.45'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +3.98%  '

# Row 49
$ws.Range('E49').Value = '  +2.15%  '

# Row 50
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0924'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +2.07%  '

# Row 51
$ws.Range('E51').Value = '  +1.87%  '
